$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2917.1396
$ws.Range("I106").Value = 2704.75
$ws.Range("J106").Value = 2999.3547
$ws.Range("K106").Value = 2704.75
$ws.Range("L106").Value = 2999.3547
$ws.Range("M106").Value = -2073.75
$ws.Range("N106").Value = -4261.3547

$ws.Range("H131").Value = 1207.8572
$ws.Range("I131").Value = 755.4545000000001
$ws.Range("J131").Value = 2866.6667
$ws.Range("K131").Value = 2266.3635
$ws.Range("L131").Value = 8600.000100000001
$ws.Range("M131").Value = 2773.6365
$ws.Range("N131").Value = -18680.0001

$ws.Range("H132").Value = 1330
$ws.Range("I132").Value = 1333.7142
$ws.Range("J132").Value = 1293.6
$ws.Range("K132").Value = 4001.1426
$ws.Range("L132").Value = 3880.8
$ws.Range("M132").Value = -1471.1426
$ws.Range("N132").Value = -8940.799999999999

$ws.Range("H137").Value = 3067.5469
$ws.Range("I137").Value = 2358.383
$ws.Range("J137").Value = 5028.1763
$ws.Range("K137").Value = 7075.148999999999
$ws.Range("L137").Value = 15084.5289
$ws.Range("M137").Value = -4525.148999999999
$ws.Range("N137").Value = -20184.5289

$ws.Range("H141").Value = 1910.766
$ws.Range("I141").Value = 1771.4359
$ws.Range("J141").Value = 2590
$ws.Range("K141").Value = 5314.307699999999
$ws.Range("L141").Value = 7770
$ws.Range("M141").Value = -134.3076999999994
$ws.Range("N141").Value = -18130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2015.9592
$ws.Range("I132").Value = 1827.421
$ws.Range("K132").Value = 5482.263
$ws.Range("M132").Value = -2952.263

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 939.0294
$ws.Range("I134").Value = 810.6786
$ws.Range("J134").Value = 1538
$ws.Range("K134").Value = 2432.0358
$ws.Range("L134").Value = 4614
$ws.Range("M134").Value = 102.9642000000003
$ws.Range("N134").Value = -9684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 41548.25
$ws.Range("J20").Value = 41548.25
$ws.Range("L20").Value = 41548.25
$ws.Range("N20").Value = -42020.25

$ws.Range("H30").Value = 41548.25
$ws.Range("J30").Value = 41548.25
$ws.Range("L30").Value = 41548.25
$ws.Range("N30").Value = -41730.25

$ws.Range("H31").Value = 20845.537
$ws.Range("I31").Value = 25899.072
$ws.Range("J31").Value = 3158.1667
$ws.Range("K31").Value = 25899.072
$ws.Range("L31").Value = 3158.1667
$ws.Range("M31").Value = -25604.072
$ws.Range("N31").Value = -3748.1667

$ws.Range("H34").Value = 20845.537
$ws.Range("I34").Value = 25899.072
$ws.Range("J34").Value = 3158.1667
$ws.Range("K34").Value = 25899.072
$ws.Range("L34").Value = 3158.1667
$ws.Range("M34").Value = -25697.072
$ws.Range("N34").Value = -3562.1667

$ws.Range("H99").Value = 1003
$ws.Range("I99").Value = 1003
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1003
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 495
$ws.Range("N99").ClearContents()

$ws.Range("H123").Value = 30528.572
$ws.Range("J123").Value = 30528.572
$ws.Range("L123").Value = 30528.572
$ws.Range("N123").Value = -40328.572

$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -24910

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 1003
$ws.Range("I126").Value = 1003
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3009
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -539
$ws.Range("N126").ClearContents()

$ws.Range("H128").Value = 41548.25
$ws.Range("J128").Value = 41548.25
$ws.Range("L128").Value = 41548.25
$ws.Range("N128").Value = -51508.25

$ws.Range("H129").Value = 45671
$ws.Range("J129").Value = 45671
$ws.Range("L129").Value = 45671
$ws.Range("N129").Value = -55671

$ws.Range("H130").Value = 35935
$ws.Range("J130").Value = 35935
$ws.Range("L130").Value = 35935
$ws.Range("N130").Value = -45975

$ws.Range("H131").Value = 32666.334
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 32666.334
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 32666.334
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -42746.334

$ws.Range("H132").Value = 1145.4459
$ws.Range("I132").Value = 826.2364
$ws.Range("K132").Value = 2478.7092
$ws.Range("M132").Value = 51.29079999999976

$ws.Range("H134").Value = 1428.921
$ws.Range("I134").Value = 1538.091
$ws.Range("J134").Value = 708.4
$ws.Range("K134").Value = 4614.272999999999
$ws.Range("L134").Value = 2125.2
$ws.Range("M134").Value = -2079.272999999999
$ws.Range("N134").Value = -7195.2

$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 50000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 50000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 4347.857
$ws.Range("J88").Value = 4347.857
$ws.Range("L88").Value = 13043.571
$ws.Range("N88").Value = -13899.571

$ws.Range("H91").Value = 4347.857
$ws.Range("J91").Value = 4347.857
$ws.Range("L91").Value = 13043.571
$ws.Range("N91").Value = -16007.571

$ws.Range("H99").Value = 2276.5833
$ws.Range("I99").Value = 1932.7778
$ws.Range("J99").Value = 3308
$ws.Range("K99").Value = 5798.3334
$ws.Range("L99").Value = 9924
$ws.Range("M99").Value = -3552.3334
$ws.Range("N99").Value = -14416

$ws.Range("H105").Value = 8471.5
$ws.Range("J105").Value = 8471.5
$ws.Range("L105").Value = 25414.5
$ws.Range("N105").Value = -30656.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1870.5435
$ws.Range("I136").Value = 1316.7797
$ws.Range("J136").Value = 2860.606
$ws.Range("K136").Value = 3950.3391
$ws.Range("L136").Value = 8581.818000000001
$ws.Range("M136").Value = -1400.3391
$ws.Range("N136").Value = -13681.818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 598.21313
$ws.Range("I132").Value = 377.87274
$ws.Range("J132").Value = 2618
$ws.Range("K132").Value = 1133.61822
$ws.Range("L132").Value = 7854
$ws.Range("M132").Value = 1396.38178
$ws.Range("N132").Value = -12914

$ws.Range("H136").Value = 332.88235
$ws.Range("I136").Value = 277.13953
$ws.Range("J136").Value = 632.5
$ws.Range("K136").Value = 831.41859
$ws.Range("L136").Value = 1897.5
$ws.Range("M136").Value = 1718.58141
$ws.Range("N136").Value = -6997.5
